# Refreshes the cryptos table (Coin/Link/Price/Volume) with the latest scrape,
# matching the scheduled "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price values look numeric (e.g. "1.110", "0.00001096") and Excel would
# otherwise silently convert them to real numbers, dropping trailing zeros or
# switching to scientific notation. A leading apostrophe keeps them literal text,
# same as the original inline-string cells.

$ws.Range("D2").Value = '28.163.87'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = '1.860.66'
$ws.Range("E3").Value = '  -0.58%  '
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").Value = "'312.64"
$ws.Range("E5").Value = '  -0.15%  '
$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("D7").Value = "'0.5093"
$ws.Range("E7").Value = '  -0.71%  '
$ws.Range("D8").Value = "'0.3900"
$ws.Range("E8").Value = '  +0.45%  '
$ws.Range("D9").Value = "'0.08265"
$ws.Range("E9").Value = '  -1.28%  '
$ws.Range("D10").Value = "'1.110"
$ws.Range("E10").Value = '  -0.26%  '
$ws.Range("D11").Value = "'41.57"
$ws.Range("E11").Value = '  -0.20%  '
$ws.Range("D12").Value = "'6.213"
$ws.Range("E12").Value = '  +0.31%  '
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").Value = "'20.21"
$ws.Range("E13").Value = '  -1.65%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.837.04'
$ws.Range("E14").Value = '  -1.55%  '
$ws.Range("D15").Value = "'7.195"
$ws.Range("E15").Value = '  -1.18%  '
$ws.Range("D16").Value = "'1.007"
$ws.Range("E16").Value = '  +0.27%  '
$ws.Range("D17").Value = "'0.00001096"
$ws.Range("E17").Value = '  -0.82%  '
$ws.Range("D18").Value = "'90.93"
$ws.Range("E18").Value = '  +0.06%  '
$ws.Range("D19").Value = "'0.06663"
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("D20").Value = "'17.60"
$ws.Range("E20").Value = '  -0.57%  '
$ws.Range("D21").Value = "'1.004"
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").Value = "'5.933"
$ws.Range("E22").Value = '  -1.56%  '
$ws.Range("D23").Value = '28.159.42'
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").Value = "'11.06"
$ws.Range("E24").Value = '  -0.44%  '
$ws.Range("D25").Value = "'2.232"
$ws.Range("E25").Value = '  -1.00%  '
$ws.Range("D26").Value = '2.047.56'
$ws.Range("E26").Value = '  -1.61%  '
$ws.Range("D27").Value = "'159.65"
$ws.Range("E27").Value = '  +0.85%  '
$ws.Range("D28").Value = "'20.47"
$ws.Range("E28").Value = '  -0.46%  '
$ws.Range("D29").Value = "'2.402"
$ws.Range("E29").Value = '  -2.94%  '
$ws.Range("D30").Value = "'125.25"
$ws.Range("E30").Value = '  +0.30%  '
$ws.Range("D31").Value = "'0.1053"
$ws.Range("E31").Value = '  -0.68%  '
$ws.Range("D32").Value = "'1.033"
$ws.Range("E32").Value = '  -0.36%  '
$ws.Range("D33").Value = "'5.813"
$ws.Range("E33").Value = '  -1.45%  '
$ws.Range("D34").Value = "'3.600"
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("D35").Value = "'0.02426"
$ws.Range("E35").Value = '  -0.26%  '
$ws.Range("D36").Value = "'0.06460"
$ws.Range("E36").Value = '  -0.96%  '
$ws.Range("D37").Value = "'9.045"
$ws.Range("E37").Value = '  -5.53%  '
$ws.Range("D38").Value = "'0.2165"
$ws.Range("E38").Value = '  -0.87%  '
$ws.Range("D39").Value = "'1.245"
$ws.Range("E39").Value = '  +1.40%  '
$ws.Range("D40").Value = "'0.6415"
$ws.Range("E40").Value = '  -1.19%  '
$ws.Range("D41").Value = "'1.180"
$ws.Range("E41").Value = '  -2.12%  '
$ws.Range("D42").Value = "'4.939"
$ws.Range("E42").Value = '  -1.27%  '
$ws.Range("D43").Value = "'11.04"
$ws.Range("E43").Value = '  -2.54%  '
$ws.Range("D44").Value = "'0.5998"
$ws.Range("E44").Value = '  -1.27%  '
$ws.Range("D45").Value = "'12.98"
$ws.Range("E45").Value = '  -1.03%  '
$ws.Range("D46").Value = "'3.658"
$ws.Range("E46").Value = '  -0.47%  '
$ws.Range("D47").Value = "'1.267"
$ws.Range("E47").Value = '  -0.79%  '
$ws.Range("D48").Value = "'1.997"
$ws.Range("E48").Value = '  -0.36%  '
$ws.Range("D49").Value = "'1.203"
$ws.Range("E49").Value = '  -1.03%  '
$ws.Range("D50").Value = "'120.68"
$ws.Range("E50").Value = '  -0.66%  '
$ws.Range("D51").Value = "'0.06867"
$ws.Range("E51").Value = '  -0.09%  '
